# v0.7.4f: Added named zombies and swords
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

# --- Insert 2 new rows right after row 8 (before the blank separator row 9) ---
$ws.Rows.Item(9).Resize(2).Insert()

# --- Append 4 new named-sword rows after the existing last row (1304, now row 38) ---
$ws.Range("B39").Value = 1305
$ws.Range("C39").Value = "Ethan Pitney"
$ws.Range("F39").Value = 5

# New row 9: unit 1007 "Michael Schmiesing" (Human, Ben)
$ws.Range("B9").Value = 1007
$ws.Range("C9").Value = "Michael Schmiesing"
$ws.Range("D9").Value = "Human"
$ws.Range("E9").Value = "Ben"
$ws.Range("F9").Value = 7
$ws.Range("H9").Value = 15
$ws.Range("I9").Value = 4
$ws.Range("J9").Value = 3
$ws.Range("L9").Value = 1.5

# New row 10: unit 1008 "Molly Schmiesing" (Human, Ben)
$ws.Range("B10").Value = 1008
$ws.Range("C10").Value = "Molly Schmiesing"
$ws.Range("D10").Value = "Human"
$ws.Range("E10").Value = "Ben"
$ws.Range("F10").Value = 6
$ws.Range("H10").Value = 9
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 1
$ws.Range("L10").Value = 1.7

$ws.Range("B40").Value = 1306
$ws.Range("C40").Value = "James Sarlo"
$ws.Range("F40").Value = 6

$ws.Range("B41").Value = 1307
$ws.Range("C41").Value = "Matt Hair"
$ws.Range("F41").Value = 7

$ws.Range("B42").Value = 1308
$ws.Range("C42").Value = "Nicholas Belt"
$ws.Range("F42").Value = 8

# Row 8 (unit 1006) gets a new drop
$ws.Range("Q8").Value = "cigar, lighter"

# Row 13 (was 1102 "Dan Gray" Frog, shifted down by the insert above) gets new stats
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2
$ws.Range("M13").Value = 3
$ws.Range("P13").Value = "agility=2"

# Column Q got wider to fit the new drop text
$ws.Columns.Item(17).ColumnWidth = 40.6

# Selection moved to H13 on this sheet
$ws.Range("H13").Select()
